$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Q0..Q8 in B1:J1
$ws.Range("B1").Value = "Q0"
$ws.Range("C1").Value = "Q1"
$ws.Range("D1").Value = "Q2"
$ws.Range("E1").Value = "Q3"
$ws.Range("F1").Value = "Q4"
$ws.Range("G1").Value = "Q5"
$ws.Range("H1").Value = "Q6"
$ws.Range("I1").Value = "Q7"
$ws.Range("J1").Value = "Q8"

# New header cell J1 needs the same style as the rest of the header row
# (bold font, centered/top aligned, thin border all around)
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").VerticalAlignment = -4160
$ws.Range("J1").Borders.LineStyle = 1

# Column A labels + numeric data per row
$ws.Range("A2").Value = "2022-01-01 00:00:00_diff"
$ws.Range("B2").Value = -2.515288930168324
$ws.Range("C2").Value = 5.469485029821115
$ws.Range("D2").Value = 3.540762996875145
$ws.Range("E2").Value = 4.959184400694483
$ws.Range("F2").Value = -1.825783515119027
$ws.Range("G2").Value = 1.874837942983618
$ws.Range("H2").Value = 1.675323243939076

$ws.Range("A3").Value = "2022-04-01 00:00:00_diff"
$ws.Range("B3").Value = 1.765433747299823
$ws.Range("C3").Value = -0.1632882856461468
$ws.Range("D3").Value = 1.255133118173191
$ws.Range("E3").Value = -5.529834797640319
$ws.Range("F3").Value = -1.829213339537674
$ws.Range("G3").Value = -2.028728038582216

$ws.Range("A4").Value = "2022-07-01 00:00:00_diff"
$ws.Range("B4").Value = 1.189993269783784
$ws.Range("C4").Value = 2.608414673603122
$ws.Range("D4").Value = -4.176553242210388
$ws.Range("E4").Value = -0.4759317841077432
$ws.Range("F4").Value = -0.6754464831522852
$ws.Range("G4").Value = 1.478677874084156
$ws.Range("H4").Value = -2.955658350684927
$ws.Range("I4").Value = 1.765133152079301
$ws.Range("J4").Value = 0.1034734828819666

$ws.Range("A5").Value = "2022-10-01 00:00:00_diff"
$ws.Range("B5").Value = 2.350441341333109
$ws.Range("C5").Value = -4.434526574480401
$ws.Range("D5").Value = -0.7339051163777555
$ws.Range("E5").Value = -0.9334198154222975
$ws.Range("F5").Value = 1.220704541814143
$ws.Range("G5").Value = -3.21363168295494
$ws.Range("H5").Value = 1.507159819809289
$ws.Range("I5").Value = -0.1544998493880456

$ws.Range("A6").Value = "2023-01-01 00:00:00_diff"
$ws.Range("B6").Value = -3.293640845203015
$ws.Range("C6").Value = 0.4069806128996306
$ws.Range("D6").Value = 0.2074659138550886
$ws.Range("E6").Value = 2.361590271091529
$ws.Range("F6").Value = -2.072745953677554
$ws.Range("G6").Value = 2.648045549086675
$ws.Range("H6").Value = 0.9863858798893403

$ws.Range("A7").Value = "2023-04-01 00:00:00_diff"
$ws.Range("B7").Value = -1.315328504146407
$ws.Range("C7").Value = -1.514843203190949
$ws.Range("D7").Value = 0.6392811540454919
$ws.Range("E7").Value = -3.795055070723591
$ws.Range("F7").Value = 0.9257364320406372
$ws.Range("G7").Value = -0.735923237156697

$ws.Range("A8").Value = "2023-07-01 00:00:00_diff"
$ws.Range("B8").Value = -1.685642723295871
$ws.Range("C8").Value = 0.46848163394057
$ws.Range("D8").Value = -3.965854590828513
$ws.Range("E8").Value = 0.7549369119357152
$ws.Range("F8").Value = -0.906722757261619
$ws.Range("G8").Value = -2.296904327241691
$ws.Range("H8").Value = -1.780776948699284
$ws.Range("I8").Value = 0.8202661914371991

$ws.Range("A9").Value = "2023-10-01 00:00:00_diff"
$ws.Range("B9").Value = 0.4964247088772
$ws.Range("C9").Value = -3.937911515891883
$ws.Range("D9").Value = 0.7828799868723453
$ws.Range("E9").Value = -0.8787796823249889
$ws.Range("F9").Value = -2.268961252305061
$ws.Range("G9").Value = -1.752833873762654
$ws.Range("H9").Value = 0.8482092663738292

$ws.Range("A10").Value = "2024-01-01 00:00:00_diff"
$ws.Range("B10").Value = -2.879911695062517
$ws.Range("C10").Value = 1.840879807701711
$ws.Range("D10").Value = 0.1792201385043768
$ws.Range("E10").Value = -1.210961431475695
$ws.Range("F10").Value = -0.6948340529332881
$ws.Range("G10").Value = 1.906209087203195

$ws.Range("A11").Value = "2024-04-01 00:00:00_diff"
$ws.Range("B11").Value = 0.6927000937718053
$ws.Range("C11").Value = -0.9689595754255289
$ws.Range("D11").Value = -2.359141145405601
$ws.Range("E11").Value = -1.843013766863194
$ws.Range("F11").Value = 0.7580293732732892

$ws.Range("A12").Value = "2024-07-01 00:00:00_diff"
$ws.Range("B12").Value = -0.2989757314280439
$ws.Range("C12").Value = -1.689157301408116
$ws.Range("D12").Value = -1.173029922865709
$ws.Range("E12").Value = 1.428013217270774

$ws.Range("A13").Value = "2024-10-01 00:00:00_diff"
$ws.Range("B13").Value = -1.407956512576945
$ws.Range("C13").Value = -0.8918291340345377
$ws.Range("D13").Value = 1.709214006101945

$ws.Range("A14").Value = "2025-01-01 00:00:00_diff"
$ws.Range("B14").Value = -0.7526570575950728
$ws.Range("C14").Value = 1.84838608254141

$ws.Range("A15").Value = "2025-04-01 00:00:00_diff"
$ws.Range("B15").Value = 2.068463645145983

$ws.Range("A16").Value = "2025-07-01 00:00:00_diff"

